# Append a new response row (row 18) to the docenti survey export sheet,
# matching the "upload docenti e studenti" data extension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that look numeric but must stay as Text (same convention used by
# every other row already in the sheet: Status, Progress, Duration,
# Finished, RecordedDate, LocationLatitude/Longitude, Q9_1/Q9_2).
$textCols = @("A","B","C","D","E","F","G","H","I","N","O","P","Q","R","S")
foreach ($col in $textCols) {
    $ws.Range("${col}18").NumberFormat = "@"
}

$ws.Range("A18").Value = "2024-12-04 10:27:51"
$ws.Range("B18").Value = "2024-12-04 10:32:19"
$ws.Range("C18").Value = "0"
$ws.Range("D18").Value = "37.159.58.54"
$ws.Range("E18").Value = "100"
$ws.Range("F18").Value = "268"
$ws.Range("G18").Value = "1"
$ws.Range("H18").Value = "1733308340.012"
$ws.Range("I18").Value = "R_225qnrdRgp00Mmt"
$ws.Range("N18").Value = "42.9786"
$ws.Range("O18").Value = "13.871"
$ws.Range("P18").Value = "anonymous"
$ws.Range("Q18").Value = "IT"
$ws.Range("R18").Value = "1"
$ws.Range("S18").Value = "1"

$ws.Range("T18").Value = 25
$ws.Range("U18").Value = 17
$ws.Range("V18").Value = 8
$ws.Range("W18").Value = 2
$ws.Range("X18").Value = 1
$ws.Range("Y18").Value = 1

$ws.Range("Z18").Value = "Buono"
$ws.Range("AA18").Value = "Buono"
$ws.Range("AB18").Value = "Buono"
$ws.Range("AC18").Value = "Buono"

$ws.Range("AD18").Value = 3

$ws.Range("AE18").Value = "Sufficiente"
$ws.Range("AF18").Value = "Sufficiente"
$ws.Range("AG18").Value = "BR04"
$ws.Range("AH18").Value = "Buono"
